# Apply crypto price/volume refresh per the GitHub Actions data update.
# Leading "'" is used on D-column values that look like plain numbers so
# Excel keeps storing them as text (matching the source inline-string
# cells) instead of silently converting them to floating point numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "97.360.83"
$ws.Range("E2").Value = "  +2.45%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.586.53"
$ws.Range("E3").Value = "  +0.59%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.07%  "

# Row 5 - Solana
$ws.Range("D5").Value = "'241.68"
$ws.Range("E5").Value = "  +2.61%  "

# Row 6 - was BNB, now XRP
$ws.Range("B6").Value = "XRP"
$ws.Range("C6").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D6").Value = "'1.73"
$ws.Range("E6").Value = "  +17.84%  "

# Row 7 - was XRP, now BNB
$ws.Range("B7").Value = "BNB"
$ws.Range("C7").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D7").Value = "'653.95"
$ws.Range("E7").Value = "  -0.22%  "

# Row 8 - Dogecoin
$ws.Range("D8").Value = "'0.433"
$ws.Range("E8").Value = "  +8.81%  "

# Row 9 - USDC
$ws.Range("E9").Value = "  -0.12%  "

# Row 10 - Cardano
$ws.Range("E10").Value = "  +5.55%  "

# Row 11 - LidoStakedEther
$ws.Range("D11").Value = "3.582.74"
$ws.Range("E11").Value = "  +0.54%  "

# Row 12 - Avalanche
$ws.Range("D12").Value = "'44.40"
$ws.Range("E12").Value = "  +5.05%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  +1.00%  "

# Row 14 - Toncoin
$ws.Range("D14").Value = "'6.47"
$ws.Range("E14").Value = "  +0.55%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "4.253.03"
$ws.Range("E15").Value = "  +0.50%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "97.107.69"
$ws.Range("E16").Value = "  +2.31%  "

# Row 17 - ShibaInu
$ws.Range("E17").Value = "  +3.48%  "

# Row 18 - Polkadot
$ws.Range("D18").Value = "'8.75"
$ws.Range("E18").Value = "  +12.97%  "

# Row 19 - WrappedEther
$ws.Range("D19").Value = "3.571.48"
$ws.Range("E19").Value = "  -0.05%  "

# Row 20 - Uniswap
$ws.Range("D20").Value = "'12.62"
$ws.Range("E20").Value = "  -0.15%  "

# Row 21 - Chainlink
$ws.Range("D21").Value = "'18.18"
$ws.Range("E21").Value = "  +2.55%  "

# Row 22 - Stellar
$ws.Range("D22").Value = "'0.531"
$ws.Range("E22").Value = "  +11.18%  "

# Row 23 - SuiNetwork
$ws.Range("D23").Value = "'3.50"
$ws.Range("E23").Value = "  +1.35%  "

# Row 24 - BitcoinCash
$ws.Range("D24").Value = "'518.27"
$ws.Range("E24").Value = "  +2.03%  "

# Row 25 - PEPE
$ws.Range("D25").Value = "'0.0000207"
$ws.Range("E25").Value = "  +6.27%  "

# Row 26 - NEARProtocol
$ws.Range("D26").Value = "'6.97"
$ws.Range("E26").Value = "  +2.59%  "

# Row 27 - Litecoin
$ws.Range("D27").Value = "'101.98"
$ws.Range("E27").Value = "  +7.35%  "

# Row 28 - Aptos
$ws.Range("D28").Value = "'13.15"
$ws.Range("E28").Value = "  +4.29%  "

# Row 29 - WrappedeETH
$ws.Range("D29").Value = "3.779.71"
$ws.Range("E29").Value = "  +0.57%  "

# Row 30 - Hedera
$ws.Range("D30").Value = "'0.172"
$ws.Range("E30").Value = "  +20.33%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  -1.11%  "

# Row 32 - InternetComputer(DFINITY)
$ws.Range("D32").Value = "'12.04"
$ws.Range("E32").Value = "  +5.00%  "

# Row 33 - Dai
$ws.Range("D33").Value = "'0.999"
$ws.Range("E33").Value = "  -0.02%  "

# Row 34 - Cronos
$ws.Range("D34").Value = "'0.186"
$ws.Range("E34").Value = "  +5.05%  "

# Row 35 - Binance-PegBSC-USD
$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = "  -0.02%  "

# Row 36 - EthereumClassic
$ws.Range("D36").Value = "'32.05"
$ws.Range("E36").Value = "  +0.73%  "

# Row 37 - was RenderToken, now PolygonEcosystemToken
$ws.Range("B37").Value = "PolygonEcosystemToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D37").Value = "'0.573"
$ws.Range("E37").Value = "  +2.98%  "

# Row 38 - Bittensor
$ws.Range("D38").Value = "'616.79"
$ws.Range("E38").Value = "  +6.80%  "

# Row 39 - was PolygonEcosystemToken, now RenderToken
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D39").Value = "'8.78"
$ws.Range("E39").Value = "  +3.63%  "

# Row 40 - Fetch.AI
$ws.Range("E40").Value = "  -1.49%  "

# Row 41 - Kaspa
$ws.Range("E41").Value = "  +2.86%  "

# Row 42 - ImmutableX
$ws.Range("D42").Value = "'1.95"
$ws.Range("E42").Value = "  +7.08%  "

# Row 43 - USDe
$ws.Range("E43").Value = "  -0.03%  "

# Row 44 - ARBITRUM
$ws.Range("D44").Value = "'0.931"
$ws.Range("E44").Value = "  +2.83%  "

# Row 45 - Filecoin
$ws.Range("D45").Value = "'6.04"
$ws.Range("E45").Value = "  +5.41%  "

# Row 46 - VeChain
$ws.Range("D46").Value = "'0.0443"
$ws.Range("E46").Value = "  +7.33%  "

# Row 47 - Stacks
$ws.Range("E47").Value = "  +2.53%  "

# Row 48 - Algorand
$ws.Range("E48").Value = "  +38.75%  "

# Row 49 - WhiteBITCoin
$ws.Range("D49").Value = "'23.64"
$ws.Range("E49").Value = "  +1.06%  "

# Row 50 - Cosmos
$ws.Range("D50").Value = "'8.57"
$ws.Range("E50").Value = "  +5.32%  "

# Row 51 - was dogwifhat, now EnergySwap
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "'33.01"
$ws.Range("E51").Value = "  -5.16%  "
